# Change table "id='...'" directives to "class='...'" directives
# across the ObjTables fixture sheets.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("!!Main root")
$sheet1.Range("A2").Value = "!!ObjTables type='Data' class='MainRoot'"

$sheet2 = $wb.Worksheets.Item("!!Nodes")
$sheet2.Range("A1").Value = "!!ObjTables type='Data' class='Node'"

$sheet3 = $wb.Worksheets.Item("!!Leaves")
$sheet3.Range("A1").Value = "!!ObjTables type='Data' class='Leaf'"

$sheet4 = $wb.Worksheets.Item("!!One to many rows")
$sheet4.Range("A1").Value = "!!ObjTables type='Data' class='OneToManyRow'"
